# Edit: Changed fraction npc dialogue
# - ratioNpc3's "How can I help" text gets a trailing "?"
# - fractionNpc2's long text is shortened/reworded
# - Two brand-new NPC dialogue rows (fractionNpc3 / fractionNpc4) are inserted
#   right after the fractionNpc2 row (55-59 range), pushing the groundskeeper
#   row (and everything below it) down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) "How can I help" -> "How can I help?" (ratioNpc3 row, column D)
$ws.Range("D54").Value = "How can I help?"

# 2) Reword fractionNpc2's instruction text (row 59, column D)
$ws.Range("D59").Value = "I need you to rescue each student that is trapped.  Use these magical rods to reach each student. `nIf you need to start again or get stuck come back to me!"

# 3) Insert two new blank rows at 60/61 (format copied from the row above,
#    matching row 59's styling), pushing the old row 60 (groundskeeper) down
#    to row 62, and every following row down by two as well.
$ws.Range("A60:A61").EntireRow.Insert() | Out-Null

# 4) Populate the new fractionNpc3 row (60)
$ws.Range("A60").Value = "fractionNpc3"
$ws.Range("B60").Value = "fractionNpc"
$ws.Range("C60").Value = "Professor of Cartography"
$ws.Range("D60").Value = "How can I help?"
$ws.Range("G60").Value = "Start Again"
$ws.Range("H60").Value = "scriptevent fraction:npc 1"
$ws.Range("I60").Value = "I’m stuck"
$ws.Range("J60").Value = "dialogue open @e[tag=fractionNpc] @p fractionNpc4"

# 5) Populate the new fractionNpc4 row (61)
$ws.Range("A61").Value = "fractionNpc4"
$ws.Range("B61").Value = "fractionNpc"
$ws.Range("C61").Value = "Professor of Cartography"
$ws.Range("D61").Value = "1. The grids are 24x24 blocks wide and you need to use the rods to cross them.`n2. You can place the rods by right clicking in front of the white blocks. `n3. Talk to each student as you go, they will let you know if you have placed the most optimum rod down. "
$ws.Range("G61").Value = "Thank you"

# Keep both new rows at the sheet's default (non-autofit) row height, same as
# the source file.
$ws.Rows.Item(60).RowHeight = 15.75
$ws.Rows.Item(61).RowHeight = 15.75

# 6) Restore view selection to roughly where the author left off editing
$ws.Range("D59").Select() | Out-Null
